# Apply the "previous gates demo" edit:
#  - survey sheet: remove the "plot_id" (integer) question row
#  - settings sheet: add a new "table_id" -> "plot" setting row
#  - update the remembered cell selections on the survey/settings sheets

$wb = $excel.ActiveWorkbook

$survey   = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# 1. Remove the plot_id question (row 2) from the survey sheet; remaining rows shift up.
[void]$survey.Rows("2:2").Delete()

# 2. Add the new table_id setting row on the settings sheet.
$settings.Range("A6").Value = "table_id"
$settings.Range("B6").Value = "plot"

# 3. Restore the recorded selections. Select settings last so it stays the active tab.
[void]$survey.Range("B7").Select()
[void]$settings.Range("B7").Select()
